$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.047.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.64%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.787.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.94%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4514'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.48%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3605'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07526'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.42'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.109'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.087'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.263'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.788.79'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.83%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.98'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.89%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001068'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06460'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.852'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.50%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.092.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.60%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.089'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.58%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.46'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.992.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.89%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.246'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.109'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09207'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.630'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.677'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02306'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.67%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06137'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2103'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6380'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.997'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.191'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.402'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.021'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5938'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.739'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '123.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.39%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.970'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.92%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.150'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06956'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.78%  '

